$p = $ppt.ActivePresentation

# The edit inserts one brand-new slide into the deck, right after the
# current slide 5 ("C# Post Roslyn") and before the old slide 6 (the C# 7
# feature-demo code slide). Every slide that used to sit at position 6-9
# simply shifts down to 7-10 — their content is untouched. The new slide
# uses the "Title and Content" layout and is left blank (empty Title +
# Content placeholders), exactly as a freshly inserted, not-yet-authored
# slide would look.

$master = $p.SlideMaster

# Find the "Title and Content" layout by name (fall back to the standard
# 2nd layout slot used by every default PowerPoint theme).
$titleAndContentLayout = $null
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $candidate = $master.CustomLayouts.Item($i)
    if ($candidate.Name -eq "Title and Content") {
        $titleAndContentLayout = $candidate
        break
    }
}
if ($null -eq $titleAndContentLayout) {
    $titleAndContentLayout = $master.CustomLayouts.Item(2)
}

$newSlide = $p.Slides.AddSlide(6, $titleAndContentLayout)
